$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 75 ("Fruta / hortaliza, semanal").
# Every existing data row from 75 down to 119 shifts down by one row (to 76..120);
# the new record's values are then written into row 75.

$srcRange = $ws.Range("A75:R119")
$destRange = $ws.Range("A76:R120")
$destRange.Value = $srcRange.Value2

$ws.Range("D75").Value = 44460
$ws.Range("J75").Value = 120
$ws.Range("K75").Value = 16000
$ws.Range("L75").Value = 16000
$ws.Range("M75").Value = 16000
$ws.Range("P75").Value = 320

# The shift-copy above only moves values, not formatting; row 120 is a brand
# new row so give its date cell the same date number-format as the rest of
# column D.
$ws.Range("D120").NumberFormat = $ws.Range("D119").NumberFormat
